$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp banner in A1
$ws.Range("A1").Value = "Datos actualizados a 4 de Octubre de 2020 a las 07:22"

# India (row 5) - refreshed totals
$ws.Range("B5").Value = 6549373
$ws.Range("C5").Value = 1960
$ws.Range("D5").Value = 5509966
$ws.Range("E5").Value = 937595

# Uzbekistan (row 59) - refreshed totals
$ws.Range("B59").Value = 58421
$ws.Range("C59").Value = 183
$ws.Range("E59").Value = 3088
$ws.Range("G59").Value = 2
$ws.Range("H59").Value = 479

# Kirguistan (row 66) - refreshed totals
$ws.Range("B66").Value = 47428
$ws.Range("C66").Value = 244
$ws.Range("D66").Value = 43418
$ws.Range("E66").Value = 2944

# Belice climbs above Uruguay / Principado de Andorra with refreshed data,
# pushing those two countries down one row each (their own figures are
# unchanged, only their row position shifts).
$ws.Range("A154").Value = "Belice"
$ws.Range("B154").Value = 2131
$ws.Range("C154").Value = 51
$ws.Range("D154").Value = 1346
$ws.Range("E154").Value = 756
$ws.Range("G154").Value = 1
$ws.Range("H154").Value = 29

$ws.Range("A155").Value = "Uruguay"
$ws.Range("B155").Value = 2122
$ws.Range("D155").Value = 1831
$ws.Range("E155").Value = 243
$ws.Range("H155").Value = 48

$ws.Range("A156").Value = "Principado de Andorra"
$ws.Range("B156").Value = 2110
$ws.Range("D156").Value = 1540
$ws.Range("E156").Value = 517
$ws.Range("H156").Value = 53
